# Update 'F' column (想去人数 / want-to-go count) values across all sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 714
$ws.Range("F8").Value = 2696
$ws.Range("F10").Value = 2077
$ws.Range("F11").Value = 846
$ws.Range("F14").Value = 6685
$ws.Range("F19").Value = 1524
$ws.Range("F21").Value = 1215
$ws.Range("F23").Value = 2590
$ws.Range("F24").Value = 1739
$ws.Range("F25").Value = 1112
$ws.Range("F26").Value = 1022
$ws.Range("F27").Value = 788
$ws.Range("F28").Value = 1121
$ws.Range("F30").Value = 5385
$ws.Range("F31").Value = 290
$ws.Range("F33").Value = 1272
$ws.Range("F35").Value = 3768
$ws.Range("F41").Value = 1058
$ws.Range("F44").Value = 915
$ws.Range("F50").Value = 90

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 399
$ws.Range("F20").Value = 607
$ws.Range("F21").Value = 261
$ws.Range("F22").Value = 359
$ws.Range("F29").Value = 69
$ws.Range("F37").Value = 210

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 405
$ws.Range("F7").Value = 1475
$ws.Range("F11").Value = 317
$ws.Range("F12").Value = 580
$ws.Range("F13").Value = 687
$ws.Range("F14").Value = 1200

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 405
$ws.Range("F9").Value = 714
$ws.Range("F10").Value = 2696
$ws.Range("F12").Value = 2077
$ws.Range("F13").Value = 846
$ws.Range("F17").Value = 6685
$ws.Range("F18").Value = 580
$ws.Range("F21").Value = 687
$ws.Range("F25").Value = 2590
$ws.Range("F26").Value = 261
$ws.Range("F27").Value = 1739
$ws.Range("F28").Value = 1112
$ws.Range("F29").Value = 1121
$ws.Range("F31").Value = 5385
$ws.Range("F32").Value = 290
$ws.Range("F34").Value = 1272
$ws.Range("F35").Value = 3768
$ws.Range("F40").Value = 69
$ws.Range("F43").Value = 915
$ws.Range("F47").Value = 210
$ws.Range("F48").Value = 210
$ws.Range("F51").Value = 90
